$d = $word.ActiveDocument

$d.Content.Find.Execute("633÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "305÷9=", 2) | Out-Null
$d.Content.Find.Execute("250÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "698÷2=", 2) | Out-Null
$d.Content.Find.Execute("407÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "498÷3=", 2) | Out-Null
$d.Content.Find.Execute("523÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "654÷4=", 2) | Out-Null
$d.Content.Find.Execute("885÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "878÷6=", 2) | Out-Null
$d.Content.Find.Execute("762÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷4=", 2) | Out-Null
$d.Content.Find.Execute("972÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "134÷9=", 2) | Out-Null
$d.Content.Find.Execute("528÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "117÷3=", 2) | Out-Null
$d.Content.Find.Execute("493÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "934÷4=", 2) | Out-Null
$d.Content.Find.Execute("134÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "655÷5=", 2) | Out-Null
$d.Content.Find.Execute("475÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "167÷3=", 2) | Out-Null
$d.Content.Find.Execute("625÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "484÷8=", 2) | Out-Null
$d.Content.Find.Execute("458÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "490÷4=", 2) | Out-Null
$d.Content.Find.Execute("984÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "302÷4=", 2) | Out-Null
$d.Content.Find.Execute("568÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "918÷3=", 2) | Out-Null
$d.Content.Find.Execute("345÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "870÷9=", 2) | Out-Null
$d.Content.Find.Execute("387÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "500÷2=", 2) | Out-Null
$d.Content.Find.Execute("900÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "791÷7=", 2) | Out-Null
$d.Content.Find.Execute("858÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "957÷4=", 2) | Out-Null
$d.Content.Find.Execute("512÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "533÷4=", 2) | Out-Null
$d.Content.Find.Execute("268÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "332÷6=", 2) | Out-Null
$d.Content.Find.Execute("296÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "288÷2=", 2) | Out-Null
$d.Content.Find.Execute("623÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "599÷7=", 2) | Out-Null
$d.Content.Find.Execute("389÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "170÷8=", 2) | Out-Null
$d.Content.Find.Execute("713÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "159÷3=", 2) | Out-Null

Write-Output "Done applying replacements"
